$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new "Player Info" sheet in front of "ODI Batting".
#    We copy the existing "ODI Batting" sheet (placing the copy right
#    before it) so the new sheet starts out with identical styles /
#    page setup, then we overwrite its contents with the player info
#    table and drop the left-over cells that came from the copy.
#
#    NOTE: after Copy()/Add() the runtime can leave old worksheet
#    variables pointing at the wrong sheet, so we always re-fetch
#    sheets by name right before using them.
# ------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Copy($battingSheet)

$playerInfo = $wb.Worksheets.Item("ODI Batting (2)")
$playerInfo.Name = "Player Info"

$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$playerInfo.Range("A2").Value = "5658"
$playerInfo.Range("B2").Value = "Saqib Mahmood"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

# remove the remaining cells/columns inherited from the "ODI Batting" copy
$playerInfo.Range("E1:J9").Clear()
$playerInfo.Range("A3:D9").Clear()

# ------------------------------------------------------------------
# 2) Rename MATCH_CARD_LINK -> MATCH_CODE and replace the full URL
#    values with the bare numeric match code, on both remaining
#    sheets.
# ------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"
$nRowsBatting = $battingSheet.UsedRange.Rows.Count
for ($row = 2; $row -le $nRowsBatting; $row++) {
    $cell = $battingSheet.Cells.Item($row, 4)
    $val = $cell.Value2
    if ($val -match "MatchCode=(\d+)") {
        $cell.Value = $matches[1]
    }
}

$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$nRowsBowling = $bowlingSheet.UsedRange.Rows.Count
for ($row = 2; $row -le $nRowsBowling; $row++) {
    $cell = $bowlingSheet.Cells.Item($row, 2)
    $val = $cell.Value2
    if ($val -match "MatchCode=(\d+)") {
        $cell.Value = $matches[1]
    }
}
